$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the aging-bucket figure for invoice row 11 (Bucket_90 -> now matches
# Bucket_60's 596141 instead of the old 12392968) and add the corresponding
# Bucket_90 entry that was missing.
$ws.Range("C11").Value = 596141
$ws.Range("E11").Value = 596141

# Move the active selection to E15, reflecting where the user left off.
$ws.Range("E15").Select() | Out-Null
